$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "x" mark from H6 down to H7, H8, H9 (connection between tiles fix)
$ws.Range("H6").Value = $null
$ws.Range("H7").Value = "x"
$ws.Range("H8").Value = "x"
$ws.Range("H9").Value = "x"

# Update the active cell selection to I6 (as reflected in the saved view state)
$ws.Range("I6").Select()
